# The document's only paragraph ends with a run that contains nothing
# but a non-breaking space (U+00A0), styled with the "apple-converted-
# space" character style - a leftover artifact from pasting text off a
# web page. That whole run is removed from the paragraph.

$d = $word.ActiveDocument

# Locate the trailing non-breaking space anywhere in the document body
# and select exactly that character/run.
$hit = $d.Content
$found = $hit.Find.Execute([char]0x00A0, $true, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)

if ($found) {
    # Deleting the matched range removes the run (<w:r>/<w:t>) that
    # holds the non-breaking space from the paragraph, just like
    # selecting that trailing space in Word and pressing Delete.
    $hit.Delete()
}
